$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (workbook.xml <sheet name=.../>)
$ws.Name = "BetaFiberA"

# 2. Small floating point precision corrections on row 13
$ws.Range("C13").Value = 0.9935434098458009
$ws.Range("D13").Value = 0.9981320030533807
$ws.Range("G13").Value = 0.9935434098458009
$ws.Range("H13").Value = 0.9981320030533807
$ws.Range("M13").Value = 0.9935154270865922
$ws.Range("O13").Value = 0.9922538338368285
$ws.Range("P13").Value = 0.9915213352877553

# 3. Small floating point precision correction on row 15
$ws.Range("K15").Value = 0.952845175639116

# 4. Add new row 16 of data (A16:P16)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 1.431032689423501
$ws.Cells.Item(16, 4).Value = 0.9794719751908326
$ws.Cells.Item(16, 5).Value = 1.214231639073742
$ws.Cells.Item(16, 6).Value = 0.8408026335015744
$ws.Cells.Item(16, 7).Value = 1.431032689423501
$ws.Cells.Item(16, 8).Value = 0.9794719751908326
$ws.Cells.Item(16, 9).Value = 1.091868016818288
$ws.Cells.Item(16, 10).Value = 0.8489512606847928
$ws.Cells.Item(16, 11).Value = 1.019012699760681
$ws.Cells.Item(16, 12).Value = 0.8639407924363697
$ws.Cells.Item(16, 13).Value = 1.431047313487744
$ws.Cells.Item(16, 14).Value = 1.096851807132287
$ws.Cells.Item(16, 15).Value = 1.116384734297413
$ws.Cells.Item(16, 16).Value = 1.036163963361223

# Copy the formatting of A15 (bold, bordered, centered) onto the new A16 cell
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
